$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 83572.5
$ws.Range("I55").Value = 111211.555
$ws.Range("J55").Value = 655.3333
$ws.Range("K55").Value = 111211.555
$ws.Range("L55").Value = 655.3333
$ws.Range("M55").Value = -110997.555
$ws.Range("N55").Value = -1083.3333
$ws.Range("H70").Value = 251312.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 251312.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 753937.5
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -754477.5
$ws.Range("H73").Value = 251312.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 251312.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 753937.5
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -755809.5
$ws.Range("H92").Value = 170.17647
$ws.Range("I92").Value = 121.666664
$ws.Range("J92").Value = 286.6
$ws.Range("K92").Value = 121.666664
$ws.Range("L92").Value = 286.6
$ws.Range("M92").Value = 1126.333336
$ws.Range("N92").Value = -2782.6
$ws.Range("H100").Value = 10911.5625
$ws.Range("I100").Value = 3539.2
$ws.Range("J100").Value = 14262.637
$ws.Range("K100").Value = 3539.2
$ws.Range("L100").Value = 14262.637
$ws.Range("M100").Value = -2998.2
$ws.Range("N100").Value = -15344.637
$ws.Range("H112").Value = 2948.6667
$ws.Range("J112").Value = 2948.6667
$ws.Range("L112").Value = 8846.000100000001
$ws.Range("N112").Value = -11062.0001
$ws.Range("H116").Value = 4399.778
$ws.Range("I116").Value = 4199.8
$ws.Range("J116").Value = 4649.75
$ws.Range("K116").Value = 4199.8
$ws.Range("L116").Value = 4649.75
$ws.Range("M116").Value = -757.8000000000002
$ws.Range("N116").Value = -11533.75
$ws.Range("H137").Value = 4535.7856
$ws.Range("I137").Value = 4520.08
$ws.Range("K137").Value = 13560.24
$ws.Range("M137").Value = -11010.24
$ws.Range("H138").Value = 5722.7354
$ws.Range("J138").Value = 8100.048
$ws.Range("L138").Value = 24300.144
$ws.Range("N138").Value = -34580.144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4299.4463
$ws.Range("I32").Value = 3968.1453
$ws.Range("K32").Value = 3968.1453
$ws.Range("M32").Value = -3681.1453
$ws.Range("H45").Value = 2979.2666
$ws.Range("I45").Value = 2136.375
$ws.Range("J45").Value = 3942.5715
$ws.Range("K45").Value = 2136.375
$ws.Range("L45").Value = 3942.5715
$ws.Range("M45").Value = -1759.375
$ws.Range("N45").Value = -4696.5715
$ws.Range("H61").Value = 2453.6924
$ws.Range("I61").Value = 1904.3636
$ws.Range("K61").Value = 1904.3636
$ws.Range("M61").Value = -1692.3636
$ws.Range("H97").Value = 1286.3334
$ws.Range("I97").Value = 1286.3334
$ws.Range("K97").Value = 1286.3334
$ws.Range("M97").Value = -790.3334
$ws.Range("H122").Value = 5698.931
$ws.Range("I122").Value = 5690
$ws.Range("K122").Value = 17070
$ws.Range("M122").Value = -14620
$ws.Range("H124").Value = 79992
$ws.Range("J124").Value = 79992
$ws.Range("L124").Value = 79992
$ws.Range("N124").Value = -89812
$ws.Range("H125").Value = 70798
$ws.Range("J125").Value = 70798
$ws.Range("L125").Value = 70798
$ws.Range("N125").Value = -80638
$ws.Range("H132").Value = 2627.6428
$ws.Range("I132").Value = 2460.1794
$ws.Range("K132").Value = 7380.5382
$ws.Range("M132").Value = -4850.5382
$ws.Range("H136").Value = 2453.6924
$ws.Range("I136").Value = 1904.3636
$ws.Range("K136").Value = 5713.0908
$ws.Range("M136").Value = -3163.0908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2091.9473
$ws.Range("I94").Value = 1316.4667
$ws.Range("K94").Value = 1316.4667
$ws.Range("M94").Value = -865.4666999999999
$ws.Range("H130").Value = 70000
$ws.Range("J130").Value = 70000
$ws.Range("L130").Value = 70000
$ws.Range("N130").Value = -80040
$ws.Range("H134").Value = 52033.953
$ws.Range("I134").Value = 4985.7856
$ws.Range("K134").Value = 14957.3568
$ws.Range("M134").Value = -12422.3568
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2125.5
$ws.Range("J16").Value = 2058
$ws.Range("L16").Value = 2058
$ws.Range("N16").Value = -2632
$ws.Range("H31").Value = 34570.395
$ws.Range("J31").Value = 37664.168
$ws.Range("L31").Value = 37664.168
$ws.Range("N31").Value = -38254.168
$ws.Range("H34").Value = 34570.395
$ws.Range("J34").Value = 37664.168
$ws.Range("L34").Value = 37664.168
$ws.Range("N34").Value = -38068.168
$ws.Range("H58").Value = 7221.129
$ws.Range("I58").Value = 6978.5
$ws.Range("K58").Value = 6978.5
$ws.Range("M58").Value = -6775.5
$ws.Range("H86").Value = 10229.25
$ws.Range("I86").Value = 9903
$ws.Range("K86").Value = 9903
$ws.Range("M86").Value = -8780
$ws.Range("H89").Value = 10229.25
$ws.Range("I89").Value = 9903
$ws.Range("K89").Value = 49515
$ws.Range("M89").Value = -43899
$ws.Range("H113").Value = 2125.5
$ws.Range("J113").Value = 2058
$ws.Range("L113").Value = 2058
$ws.Range("N113").Value = -6398
$ws.Range("H134").Value = 480034.56
$ws.Range("I134").Value = 4229.5
$ws.Range("K134").Value = 12688.5
$ws.Range("M134").Value = -10153.5
$ws.Range("H136").Value = 7221.129
$ws.Range("I136").Value = 6978.5
$ws.Range("K136").Value = 20935.5
$ws.Range("M136").Value = -18385.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 1500
$ws.Range("I110").Value = 1500
$ws.Range("K110").Value = 4500
$ws.Range("M110").Value = -410
$ws.Range("H112").Value = 337599.66
$ws.Range("I112").Value = 999999
$ws.Range("K112").Value = 2999997
$ws.Range("M112").Value = -2998889
$ws.Range("H118").Value = 1338.5
$ws.Range("I118").Value = 1234.8334
$ws.Range("J118").Value = 1649.5
$ws.Range("K118").Value = 3704.5002
$ws.Range("L118").Value = 4948.5
$ws.Range("M118").Value = -2461.5002
$ws.Range("N118").Value = -7434.5
$ws.Range("H134").Value = 3362.7144
$ws.Range("I134").Value = 2708.8
$ws.Range("J134").Value = 4997.5
$ws.Range("K134").Value = 8126.400000000001
$ws.Range("L134").Value = 14992.5
$ws.Range("M134").Value = -3056.400000000001
$ws.Range("N134").Value = -25132.5
$ws.Range("H140").Value = 1627.7391
$ws.Range("I140").Value = 1383.5454
$ws.Range("K140").Value = 4150.6362
$ws.Range("M140").Value = 1029.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 38113.332
$ws.Range("I132").Value = 4733.524
$ws.Range("J132").Value = 115999.555
$ws.Range("K132").Value = 14200.572
$ws.Range("L132").Value = 347998.665
$ws.Range("M132").Value = -11670.572
$ws.Range("N132").Value = -353058.665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4094.95
$ws.Range("I46").Value = 3646.077
$ws.Range("J46").Value = 4928.5713
$ws.Range("K46").Value = 3646.077
$ws.Range("L46").Value = 4928.5713
$ws.Range("M46").Value = -3458.077
$ws.Range("N46").Value = -5304.5713
$ws.Range("H68").Value = 3016.6667
$ws.Range("I68").Value = 2550
$ws.Range("K68").Value = 2550
$ws.Range("M68").Value = -1801
$ws.Range("H71").Value = 3016.6667
$ws.Range("I71").Value = 2550
$ws.Range("K71").Value = 12750
$ws.Range("M71").Value = -9006
$ws.Range("H93").Value = 47620500
$ws.Range("I93").Value = 62501316
$ws.Range("J93").Value = 1894.8
$ws.Range("K93").Value = 62501316
$ws.Range("L93").Value = 1894.8
$ws.Range("M93").Value = -62500068
$ws.Range("N93").Value = -4390.8
$ws.Range("H127").Value = 88493.25
$ws.Range("J127").Value = 88493.25
$ws.Range("L127").Value = 88493.25
$ws.Range("N127").Value = -98413.25
$ws.Range("H132").Value = 5942.48
$ws.Range("I132").Value = 4835.222
$ws.Range("K132").Value = 14505.666
$ws.Range("M132").Value = -11975.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14999.5
$ws.Range("I45").Value = 14000
$ws.Range("J45").Value = 15332.667
$ws.Range("K45").Value = 14000
$ws.Range("L45").Value = 15332.667
$ws.Range("M45").Value = -13509
$ws.Range("N45").Value = -16314.667
$ws.Range("H62").Value = 7500.125
$ws.Range("I62").Value = 6667
$ws.Range("K62").Value = 6667
$ws.Range("M62").Value = -6043
$ws.Range("H65").Value = 7500.125
$ws.Range("I65").Value = 6667
$ws.Range("K65").Value = 33335
$ws.Range("M65").Value = -30215
$ws.Range("H81").Value = 2424.3
$ws.Range("I81").Value = 2138.111
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 4276.222
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -3215.222
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 2424.3
$ws.Range("I84").Value = 2138.111
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 21381.11
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -16077.11
$ws.Range("N84").Value = -60608
$ws.Range("H96").Value = 145725
$ws.Range("J96").Value = 4100
$ws.Range("L96").Value = 4100
$ws.Range("N96").Value = -6846
$ws.Range("H122").Value = 47621270
$ws.Range("I122").Value = 66668450
$ws.Range("K122").Value = 200005350
$ws.Range("M122").Value = -200002900
$ws.Range("H136").Value = 17289530
$ws.Range("I136").Value = 21486286
$ws.Range("J136").Value = 502499.25
$ws.Range("K136").Value = 64458858
$ws.Range("L136").Value = 1507497.75
$ws.Range("M136").Value = -64456308
$ws.Range("N136").Value = -1512597.75

Write-Output "Applied 252 cell updates across 8 sheets"